$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (Ano 2025) with refreshed "ADD" figures
$ws.Range("B7").Value = 2389782.08
$ws.Range("C7").Value = -46.21336795221598
$ws.Range("D7").Value = 2454
$ws.Range("E7").Value = 2454
$ws.Range("F7").Value = 973.8313284433578
$ws.Range("G7").Value = 3.803377904769834
